$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumen")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Row 20: "Número de transacciones" - updated counts for the new month
$ws.Range("B20").Value = 19461330
$ws.Range("C20").Value = 22940872
$ws.Range("D20").Value = 42402202

# Row 21: "Valor transacciones" - previously stored as text strings with
# thousands separators; now real numbers formatted with a thousands-separator
# number format (built-in format #,##0).
$ws.Range("B21").Value = 7087489431424
$ws.Range("C21").Value = 8023170364289.7402
$ws.Range("D21").Value = 15110659795713.699
$ws.Range("B21:D21").NumberFormat = "#,##0"

# Approximate the "best fit" column widths recorded for columns B:D after
# the update (values widened to fit the bigger formatted numbers).
$ws.Columns.Item(2).ColumnWidth = 15.6
$ws.Columns.Item(3).ColumnWidth = 15.6
$ws.Columns.Item(4).ColumnWidth = 16.6

# Reflect the final selection left behind on the sheet.
$ws.Range("D21").Select() | Out-Null
